$p = $ppt.ActivePresentation

# Slide 8 ("Data" title slide): re-center the title textbox horizontally on the slide.
$s8 = $p.Slides.Item(8)
$shp60 = $s8.Shapes.Item(1)
$shp60.Left = 24.543385826771654

# Slide 9 (MovieLens overview slide): group all the top-level shapes together
# and re-center the resulting group horizontally on the slide.
$s9 = $p.Slides.Item(9)
$range = $s9.Shapes.Range(@(1, 2, 3, 4, 5, 6, 7))
$grp = $range.Group()
$grp.Left = 24.543385826771654
